$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add new bibliography entry in row 8 (Pike-Wilson & Karayiannis, 2014)
$ws.Range("B8").Value = "Flow boiling of R245fa in 1.1 mm diameter stainless steel, `nbrass and copper tubes"
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Value = "[E.A. Pike-Wilson, T.G. Karayiannis]"
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = "R245fa"
$ws.Range("L8").Value = "SS"

# Update the heat-flux header unit from [w/m2/s] to [kW/m2]
$ws.Range("Q1").Value = "q""`n[kW/m²]"

$ws.Range("Q8").Value = "[10 ; 60]"
